$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing OT values on existing rows 102 and 103 ---
# These look numeric but must stay stored as text (matching the rest of
# column E / the source data), so we use the leading-apostrophe text-entry
# trick and then reset the style back to Normal so no stray number-format
# style is introduced.
$ws.Range('E102').Value = "'810712730"
$ws.Range('E102').Style = "Normal"

$ws.Range('E103').Value = "'810712782"
$ws.Range('E103').Style = "Normal"

# --- Append the new record as row 104 ---
$ws.Range('A104').Value = "'7834"
$ws.Range('A104').Style = "Normal"

$ws.Range('B104').Value = "'11/11/2025"
$ws.Range('B104').Style = "Normal"

$ws.Range('C104').Value = "MUÑECAS 1277"

$ws.Range('D104').Value = "'15"
$ws.Range('D104').Style = "Normal"

$ws.Range('E104').Value = "'810713028"
$ws.Range('E104').Style = "Normal"

$ws.Range('F104').Value = "PEBCOM"
$ws.Range('G104').Value = "Pendiente"
$ws.Range('H104').Value = "Picada"
$ws.Range('I104').Value = 1
$ws.Range('J104').Value = "Cambio"
$ws.Range('K104').Value = "Sin equipos"
$ws.Range('L104').Value = "Pasante"
$ws.Range('M104').Value = -58.44993
$ws.Range('N104').Value = -34.596737
$ws.Range('O104').Value = "Paternal"
$ws.Range('P104').Value = "Capital Norte"
$ws.Range('Q104').Value = "VCR-?"
$ws.Range('R104').Value = "Fuera de Poligono OVL"
